# A64 Security Compliance Tracker - update
# - Detailed Controls: flip several NON-COMPLIANT rows to PARTIAL/COMPLIANT with new evidence text
# - Remediation Roadmap: update statuses/notes to reflect progress
# - Evidence Log: append 5 new evidence rows (41-45)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Detailed Controls
# ---------------------------------------------------------------------------
$wsControls = $wb.Worksheets.Item("Detailed Controls")

$wsControls.Range("D27").Value = "PARTIAL"
$wsControls.Range("E27").Value = "Backup encryption implemented (AES-256-CBC with PBKDF2 100k iterations) in scripts/backup/mongodb_backup.sh. Enabled via ENCRYPT_BACKUPS=true in production docker-compose.prod.yml. Restore script updated to handle encrypted backups."

$wsControls.Range("D30").Value = "COMPLIANT"
$wsControls.Range("E30").Value = "Data Classification Policy created with 4-tier system (Docs/2-Working-Progress/Data-Classification-Policy.md). Complete PII inventory across 25+ collections."

$wsControls.Range("D35").Value = "PARTIAL"
$wsControls.Range("E35").Value = "Centralized structured JSON logging implemented (src/core/logging_config.py). JSON format in production, text in dev. Docker json-file driver with 10MB rotation. Log volumes mounted for persistence."

$wsControls.Range("D67").Value = "PARTIAL"
$wsControls.Range("E67").Value = "GitHub Actions CI/CD: security-scan.yml runs pip-audit, npm audit, Trivy container scanning on push/PR"

$wsControls.Range("D69").Value = "PARTIAL"
$wsControls.Range("E69").Value = "Reduced from 19 CVEs to 2 (protobuf pinned by google-cloud-aiplatform, ecdsa upstream wontfix). npm: 8 moderate dev-only vulnerabilities in build tools."

$wsControls.Range("D70").Value = "COMPLIANT"
$wsControls.Range("E70").Value = "GitHub Actions: security-scan.yml (lint + security scans) and build.yml (API tests, frontend build, compose validation)"

$wsControls.Range("D89").Value = "PARTIAL"
$wsControls.Range("E89").Value = "MongoDB and Redis ports removed in production (docker-compose.prod.yml: ports: []). Adminer disabled. Registry localhost-only. Ports still exposed in development."

$wsControls.Range("D100").Value = "COMPLIANT"
$wsControls.Range("E100").Value = "Data Flow Diagrams created (Docs/2-Working-Progress/Data-Flow-Diagrams.md) with system context, 5 detailed flow diagrams, data-at-rest/transit maps, retention policy."

$wsControls.Range("D104").Value = "COMPLIANT"
$wsControls.Range("E104").Value = "Security zone definitions created (Docs/2-Working-Progress/Security-Zone-Definitions.md) with 5 zones (Internet/DMZ/Application/Data/Management), traffic flow matrix, data classification boundaries."

$wsControls.Range("D105").Value = "COMPLIANT"
$wsControls.Range("E105").Value = "STRIDE threat model created (Docs/2-Working-Progress/Threat-Model.md) with 26 identified threats, attack trees, risk matrix, prioritized remediation."

# ---------------------------------------------------------------------------
# Sheet: Remediation Roadmap
# ---------------------------------------------------------------------------
$wsRoadmap = $wb.Worksheets.Item("Remediation Roadmap")

$wsRoadmap.Range("K2").Value = "Backup encryption added. Classification policy complete. Remaining: backup verification testing."

$wsRoadmap.Range("H3").Value = "In Progress"
$wsRoadmap.Range("K3").Value = "Backup encryption added. Classification policy complete. Remaining: backup verification testing."

$wsRoadmap.Range("K4").Value = "Backup encryption added. Classification policy complete. Remaining: backup verification testing."

$wsRoadmap.Range("H8").Value = "In Progress"
$wsRoadmap.Range("K8").Value = "CI/CD created, CVEs reduced. Remaining: penetration testing, WAF-based scanning."

$wsRoadmap.Range("H9").Value = "In Progress"
$wsRoadmap.Range("K9").Value = "CI/CD created, CVEs reduced. Remaining: penetration testing, WAF-based scanning."

$wsRoadmap.Range("H17").Value = "In Progress"
$wsRoadmap.Range("K17").Value = "Structured JSON logging implemented. Remaining: SIEM integration, centralized log aggregation."

$wsRoadmap.Range("H18").Value = "In Progress"
$wsRoadmap.Range("K18").Value = "Structured JSON logging implemented. Remaining: SIEM integration, centralized log aggregation."

$wsRoadmap.Range("H28").Value = "In Progress"
$wsRoadmap.Range("K28").Value = "CI/CD created, CVEs reduced. Remaining: penetration testing, WAF-based scanning."

$wsRoadmap.Range("H34").Value = "In Progress"
$wsRoadmap.Range("K34").Value = "Threat model and security zones documented. Remaining: cloud security review."

$wsRoadmap.Range("H38").Value = "In Progress"
$wsRoadmap.Range("K38").Value = "Threat model and security zones documented. Remaining: cloud security review."

# ---------------------------------------------------------------------------
# Sheet: Evidence Log - append rows 41-45
# ---------------------------------------------------------------------------
$wsEvidence = $wb.Worksheets.Item("Evidence Log")

# Note: B and H columns hold text that LOOKS like a number/date ("6", "25",
# "2026-02-05"). A leading apostrophe forces Excel to store these as literal
# text (matching the source workbook's inlineStr cells) instead of silently
# coercing them to a number/date serial; re-applying the "Normal" style
# afterwards drops the quote-prefix formatting flag so the cell is left with
# plain default formatting (same as every other untouched cell on this sheet).

$wsEvidence.Range("A41").Value = 40
$wsEvidence.Range("B41").Value = "'6"
$wsEvidence.Range("B41").Style = "Normal"
$wsEvidence.Range("C41").Value = "Data Protection"
$wsEvidence.Range("D41").Value = "Code"
$wsEvidence.Range("E41").Value = "scripts/backup/mongodb_backup.sh, scripts/backup/mongodb_restore.sh, docker-compose.prod.yml"
$wsEvidence.Range("F41").Value = "Implemented AES-256-CBC backup encryption with PBKDF2 key derivation (100k iterations) in mongodb_backup.sh. Restore script updated for encrypted backup handling. Production docker-compose enables encryption by default."
$wsEvidence.Range("G41").Value = "Claude Code"
$wsEvidence.Range("H41").Value = "'2026-02-05"
$wsEvidence.Range("H41").Style = "Normal"

$wsEvidence.Range("A42").Value = 41
$wsEvidence.Range("B42").Value = "'25"
$wsEvidence.Range("B42").Style = "Normal"
$wsEvidence.Range("C42").Value = "Security Architecture"
$wsEvidence.Range("D42").Value = "Documentation"
$wsEvidence.Range("E42").Value = "Docs/2-Working-Progress/Threat-Model.md"
$wsEvidence.Range("F42").Value = "Created STRIDE threat model with 26 threats across 6 categories, 2 attack trees, risk matrix (4 CRITICAL, 7 HIGH, 10 MEDIUM, 5 LOW), and prioritized remediation plan"
$wsEvidence.Range("G42").Value = "Claude Code"
$wsEvidence.Range("H42").Value = "'2026-02-05"
$wsEvidence.Range("H42").Style = "Normal"

$wsEvidence.Range("A43").Value = 42
$wsEvidence.Range("B43").Value = "'25"
$wsEvidence.Range("B43").Style = "Normal"
$wsEvidence.Range("C43").Value = "Security Architecture"
$wsEvidence.Range("D43").Value = "Documentation"
$wsEvidence.Range("E43").Value = "Docs/2-Working-Progress/Security-Zone-Definitions.md"
$wsEvidence.Range("F43").Value = "Created security zone definitions with 5 zones (Internet/DMZ/Application/Data/Management), trust boundaries, allowed/denied traffic matrix, data classification per zone"
$wsEvidence.Range("G43").Value = "Claude Code"
$wsEvidence.Range("H43").Value = "'2026-02-05"
$wsEvidence.Range("H43").Style = "Normal"

$wsEvidence.Range("A44").Value = 43
$wsEvidence.Range("B44").Value = "'7"
$wsEvidence.Range("B44").Style = "Normal"
$wsEvidence.Range("C44").Value = "Log Monitoring"
$wsEvidence.Range("D44").Value = "Code"
$wsEvidence.Range("E44").Value = "src/core/logging_config.py, src/main.py"
$wsEvidence.Range("F44").Value = "Implemented centralized structured JSON logging. Production uses JSON format for log aggregation. Development uses human-readable text. Respects LOG_LEVEL from environment."
$wsEvidence.Range("G44").Value = "Claude Code"
$wsEvidence.Range("H44").Value = "'2026-02-05"
$wsEvidence.Range("H44").Style = "Normal"

$wsEvidence.Range("A45").Value = 44
$wsEvidence.Range("B45").Value = "'14"
$wsEvidence.Range("B45").Style = "Normal"
$wsEvidence.Range("C45").Value = "Security Assessment"
$wsEvidence.Range("D45").Value = "Code"
$wsEvidence.Range("E45").Value = ".github/workflows/security-scan.yml, .github/workflows/build.yml"
$wsEvidence.Range("F45").Value = "GitHub Actions CI/CD pipelines created: security-scan.yml (black, flake8, pip-audit, npm audit, Trivy) and build.yml (API tests, frontend build, compose validation). Python CVEs reduced from 19 to 2."
$wsEvidence.Range("G45").Value = "Claude Code"
$wsEvidence.Range("H45").Value = "'2026-02-05"
$wsEvidence.Range("H45").Style = "Normal"
